{"js": "// Locate the (single) table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// --- Step 1: rewrite the first three single-value rows ------------------\n// Row 0: \"100\" -> \"0M\"\n// Row 1: \"0\"   -> \"0M\"\n// Row 2: \"221\" -> \"0M\"\nfor (const idx of [0, 1, 2]) {\n  const row = rows.items[idx];\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nrows.items[0].cells.items[0].value = \"0M\";\nrows.items[1].cells.items[0].value = \"0M\";\nrows.items[2].cells.items[0].value = \"0M\";\nawait context.sync();\n\n// --- Step 2: insert 10 new single-value rows right after row 2 ----------\nconst newRowValues = [\n  [\"105\"],\n  [\"0.00002\"],\n  [\"0.00005\"],\n  [\"0.00004\"],\n  [\"0.00001\"],\n  [\"0.00003\"],\n  [\"0.00004\"],\n  [\"0.00005\"],\n  [\"0.00387\"],\n  [\"100.0\"],\n];\nrows.items[2].insertRows(\"After\", newRowValues.length, newRowValues);\nawait context.sync();\n\n// --- Step 3: simplify the three multi-value (tab separated) rows --------\n// They now sit 10 rows further down than before the insert.\n// Originally at indices 33, 34, 35; now at 43, 44, 45.\nconst tables2 = context.document.body.tables;\ntables2.load(\"items\");\nawait context.sync();\nconst table2 = tables2.items[0];\nconst rows2 = table2.rows;\nrows2.load(\"items\");\nawait context.sync();\n\nfor (const idx of [43, 44, 45]) {\n  rows2.items[idx].cells.load(\"items\");\n}\nawait context.sync();\n\nrows2.items[43].cells.items[0].value = \"100\";\nrows2.items[44].cells.items[0].value = \"0\";\nrows2.items[45].cells.items[0].value = \"221\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Step 1: rewrite the first three single-value rows -------------------\n# Row 1: \"100\" -> \"0M\"\n# Row 2: \"0\"   -> \"0M\"\n# Row 3: \"221\" -> \"0M\"\n$t.Rows.Item(1).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(2).Cells.Item(1).Range.Text = \"0M\"\n$t.Rows.Item(3).Cells.Item(1).Range.Text = \"0M\"\n\n# --- Step 2: insert 10 new single-value rows right after row 3 -----------\n$newValues = @(\"105\", \"0.00002\", \"0.00005\", \"0.00004\", \"0.00001\", \"0.00003\", \"0.00004\", \"0.00005\", \"0.00387\", \"100.0\")\n\n$afterRow = $t.Rows.Item(4)\nforeach ($val in $newValues) {\n    $newRow = $t.Rows.Add($afterRow)\n    $newRow.Cells.Item(1).Range.Text = $val\n}\n\n# --- Step 3: simplify the three multi-value (tab separated) rows ---------\n# They now sit 10 rows further down than before the insert.\n# Originally rows 34/35/36 (1-based); now rows 44/45/46.\n$t.Rows.Item(44).Cells.Item(1).Range.Text = \"100\"\n$t.Rows.Item(45).Cells.Item(1).Range.Text = \"0\"\n$t.Rows.Item(46).Cells.Item(1).Range.Text = \"221\"\n"}
